$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 2 (shifts existing data rows 2-22 down to 3-23)
$ws.Rows("2:2").Insert()

# Populate the newly inserted row with the new accelerometer sample
$ws.Cells.Item(2, 1).Value = -0.647717118263246
$ws.Cells.Item(2, 2).Value = 0.6091025024652482
$ws.Cells.Item(2, 3).Value = -1.090710066258908

# The last two original data rows (now shifted to rows 22 and 23) are no
# longer part of the dataset, so remove them, leaving data in rows 2-21.
$ws.Rows("22:23").Delete()
